$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header figures -------------------------------------------------
# Salario Basico
$ws.Range("E11").Value = 690000
# Cant. Trabajadores (now only 1 worker in this batch)
$ws.Range("C13").Value = 1
# Cant. Periodos (now 12 periods, 2409..2508)
$ws.Range("F13").Value = 12

# --- Detail rows: update periods (now ascending 2409..2508) ---------
$periods = @("2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("E" + $r).Value = $periods[$i]
}

# Mora value only differs on the first period row (now 2409) and the
# last of the old block (now 2507 at row 26, back to the normal amount)
$ws.Range("F16").Value = 30000
$ws.Range("F26").Value = 60000

# --- Remove the three other workers (EDUARDO, KATHERINE, LUIS) ------
# keeping the last (closing-style) row, which shifts up to row 27 and
# gets re-purposed for JUAN DAVID VALENCIA ACOSTA's new period 2508.
$ws.Rows("27:29").Delete()

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1143365690"
$ws.Range("D27").Value = "JUAN DAVID VALENCIA ACOSTA"
$ws.Range("E27").Value = "2508"
$ws.Range("F27").Value = 60000
$ws.Range("G27").Value = 1500000

# --- Column D is now a touch narrower since the widest names were
# removed from the worker list -> let Excel recompute the best fit.
$ws.Columns("D").AutoFit()
